# ValueSet-med-thiotepa-vs.xlsx — metadata refresh (version/date/status/contacts/jurisdiction)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Simple in-place value updates (rows 1-10 keep their meaning) ---
$ws1.Range("B3").Value  = "0.1.7"
$ws1.Range("B6").Value  = "draft"
$ws1.Range("B8").Value  = "2024-11-22T12:33:30-06:00"
$ws1.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- Make room for a new "Jurisdiction" row at 12 ---
# Cascade rows 12..15 down to 13..16 (value + style) before overwriting row 12,
# so the existing style (s="2") is preserved instead of Excel minting a fresh one.
# (Copy() only overwrites a destination cell when the source has content, so
# blank source cells need an explicit ClearContents afterwards.)
$ws1.Range("A15:B15").Copy($ws1.Range("A16:B16"))

$ws1.Range("A14:B14").Copy($ws1.Range("A15:B15"))
$ws1.Range("B15").ClearContents()

$ws1.Range("A13:B13").Copy($ws1.Range("A14:B14"))
$ws1.Range("B14").ClearContents()

$ws1.Range("A12:B12").Copy($ws1.Range("A13:B13"))
$excel.CutCopyMode = $false

# Row 11 (second Contact entry)
$ws1.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# Row 12 (new Jurisdiction row)
$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = ""
